# Talent.xlsx: unify the conception of DataNode, DataTable, Entity.
# Rename the sheet and tidy up its row/column sizing + selection to match
# the re-saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to be called "Property1" - rename it to "DataNode" now
# that Property/DataTable/Entity have been folded into one concept.
$ws.Name = "DataNode"

# Header row and the second header-ish row (row 8) both got a little
# shorter (28 -> 27 pts; row 1 gained an explicit height of 27 too).
$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(8).RowHeight = 27

# Columns A and H were nudged a touch narrower.
$ws.Columns.Item(1).ColumnWidth = 20.142857142857142
$ws.Columns.Item(8).ColumnWidth = 25.428571428571427

# Leave the selection where the author left it when they saved.
[void]$ws.Range("D22").Select()
